$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 67
$ws.Range("E2").Value = 538.340118597339
$ws.Range("F2").Value = 110502
$ws.Range("G2").Value = 205.264285871757
$ws.Range("I2").Value = 0.4
$ws.Range("J2").Value = 82321
$ws.Range("K2").Value = 19053
$ws.Range("L2").Value = 1167
$ws.Range("M2").Value = 671
$ws.Range("N2").Value = 10052
$ws.Range("O2").Value = 54.9462878787879
$ws.Range("P2").Value = 6.03990958919793
$ws.Range("Q2").Value = 11.2098672384383
$ws.Range("R2").Value = 36.1194029850746
$ws.Range("S2").Value = 0.453731343283582
$ws.Range("D5").Value = 244
$ws.Range("E5").Value = 650.400166711181
$ws.Range("F5").Value = 433749
$ws.Range("G5").Value = 666.895585518219
$ws.Range("I5").Value = 0.09
$ws.Range("J5").Value = 219381
$ws.Range("K5").Value = 40602
$ws.Range("L5").Value = 3048
$ws.Range("M5").Value = 77119
$ws.Range("N5").Value = 196413
$ws.Range("O5").Value = 56.7661715481172
$ws.Range("P5").Value = 7.98943851642652
$ws.Range("Q5").Value = 10.9229564335092
$ws.Range("S5").Value = 0.45655737704918
$ws.Range("B6").Value = "East Chicago"
$ws.Range("D6").Value = 595
$ws.Range("E6").Value = 234.787367915844
$ws.Range("F6").Value = 651778
$ws.Range("G6").Value = 2776.03520915835
$ws.Range("J6").Value = 214835
$ws.Range("K6").Value = 368432
$ws.Range("L6").Value = 1437
$ws.Range("M6").Value = 5310
$ws.Range("N6").Value = 127797
$ws.Range("O6").Value = 47.5288342440801
$ws.Range("P6").Value = 10.6400467700909
$ws.Range("Q6").Value = 11.544172148924
$ws.Range("R6").Value = 29.5294117647059
$ws.Range("S6").Value = 0.369411764705882
$ws.Range("B8").Value = "Sergeant Bluff"
$ws.Range("D8").Value = 48
$ws.Range("E8").Value = 1045.26318446616
$ws.Range("F8").Value = 55152
$ws.Range("G8").Value = 52.7637448822683
$ws.Range("I8").Value = 0.38
$ws.Range("J8").Value = 43419
$ws.Range("K8").Value = 1650
$ws.Range("L8").Value = 3977
$ws.Range("M8").Value = 947
$ws.Range("N8").Value = 9775
$ws.Range("O8").Value = 65.6720425531915
$ws.Range("P8").Value = 6.48747976102618
$ws.Range("Q8").Value = 6.49071037566484
$ws.Range("S8").Value = 0.214893617021277
$ws.Range("B9").Value = "Geismar"
$ws.Range("D10").Value = 69
$ws.Range("E10").Value = 555.22662309248
$ws.Range("F10").Value = 182469
$ws.Range("G10").Value = 328.638779934022
$ws.Range("I10").Value = 0.52
$ws.Range("J10").Value = 119600
$ws.Range("K10").Value = 49943
$ws.Range("M10").Value = 3623
$ws.Range("N10").Value = 9220
$ws.Range("O10").Value = 72.2608507462687
$ws.Range("P10").Value = 5.92115845564032
$ws.Range("Q10").Value = 6.6311776732443
$ws.Range("R10").Value = 79.4202898550725
$ws.Range("S10").Value = 0.521739130434783
$ws.Range("D11").Value = 34
$ws.Range("E11").Value = 761.435287872035
$ws.Range("F11").Value = 37756
$ws.Range("G11").Value = 49.5853037038982
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0.71
$ws.Range("J11").Value = 24324
$ws.Range("K11").Value = 12493
$ws.Range("L11").Value = 145
$ws.Range("M11").Value = 136
$ws.Range("N11").Value = 653
$ws.Range("O11").Value = 34.7644117647059
$ws.Range("P11").Value = 8.12780482947082
$ws.Range("Q11").Value = 20.6481178232457
$ws.Range("D12").Value = 108
$ws.Range("E12").Value = 1195.70892071902
$ws.Range("F12").Value = 128891
$ws.Range("G12").Value = 107.794629417412
$ws.Range("I12").Value = 0.11
$ws.Range("J12").Value = 88781
$ws.Range("K12").Value = 28588
$ws.Range("L12").Value = 487
$ws.Range("N12").Value = 31670
$ws.Range("O12").Value = 55.460206185567
$ws.Range("P12").Value = 8.66775464619718
$ws.Range("Q12").Value = 8.88279091282656
$ws.Range("R12").Value = 57.962962962963
$ws.Range("S12").Value = 0.399074074074074
$ws.Range("D13").Value = 33
$ws.Range("E13").Value = 429.942996016745
$ws.Range("F13").Value = 51362
$ws.Range("G13").Value = 119.462348441186
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = 0.7
$ws.Range("J13").Value = 46644
$ws.Range("K13").Value = 1062
$ws.Range("M13").Value = 623
$ws.Range("N13").Value = 20610
$ws.Range("O13").Value = 62.32375
$ws.Range("P13").Value = 5.87796124304657
$ws.Range("Q13").Value = 5.93507974024599
$ws.Range("S13").Value = 0.212121212121212
$ws.Range("D15").Value = 70
$ws.Range("E15").Value = 584.977037050589
$ws.Range("F15").Value = 76577
$ws.Range("G15").Value = 130.905993141364
$ws.Range("I15").Value = 0.63
$ws.Range("J15").Value = 68884
$ws.Range("K15").Value = 4726
$ws.Range("M15").Value = 854
$ws.Range("N15").Value = 680
$ws.Range("O15").Value = 49.6100895522388
$ws.Range("P15").Value = 8.84973276756062
$ws.Range("Q15").Value = 9.73087658338039
$ws.Range("R15").Value = 30.8571428571429
$ws.Range("S15").Value = 0.395714285714286
